# A new price-report observation was inserted into the daily log at
# row 321 (the data is kept in reverse-ish / per-date order within the
# sheet), pushing every existing record from row 321 down to row 322,
# and so on through the former last row (444) which becomes row 445.
#
# Concretely: insert one blank row at position 321 (Excel shifts
# everything below it down by one automatically, including the used
# range / dimension), then populate that new row with the new
# observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 321:444 down to 322:445, leaving a blank row 321.
$ws.Rows(321).Insert()

# Fill in the new row 321 with the inserted observation.
$ws.Cells.Item(321, 1).Value()  = 9
$ws.Cells.Item(321, 2).Value()  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(321, 3).Value()  = "Metropolitana"
$ws.Cells.Item(321, 4).Value()  = 44795
$ws.Cells.Item(321, 5).Value()  = 13
$ws.Cells.Item(321, 6).Value()  = 100112039
$ws.Cells.Item(321, 7).Value()  = "Ciboulette"
$ws.Cells.Item(321, 8).Value()  = "Sin especificar"
$ws.Cells.Item(321, 9).Value()  = "Primera"
$ws.Cells.Item(321, 10).Value() = 250
$ws.Cells.Item(321, 11).Value() = 2000
$ws.Cells.Item(321, 12).Value() = 2000
$ws.Cells.Item(321, 13).Value() = 2000
$ws.Cells.Item(321, 14).Value() = "`$/docena de atados"
$ws.Cells.Item(321, 15).Value() = "Región Metropolitana"
$ws.Cells.Item(321, 16).Value() = 667
$ws.Cells.Item(321, 17).Value() = 3
$ws.Cells.Item(321, 18).Value() = "Hortaliza"
